# repull data, push all data, mean calculation
# Updates the "dSF" (column F) values for several rows on Sheet1 to reflect
# the re-pulled data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "F2"  = -1
    "F5"  = 2
    "F7"  = -2
    "F10" = 2
    "F15" = 4
    "F17" = -2
    "F21" = 3
    "F22" = 4
    "F26" = 2
    "F28" = 3
    "F35" = -3
    "F36" = 1
    "F46" = -2
    "F48" = 0
    "F49" = 2
    "F50" = -1
    "F52" = 0
    "F58" = 2
    "F61" = 2
    "F62" = 2
    "F63" = 1
    "F65" = 3
    "F66" = 1
    "F73" = 4
    "F74" = -2
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
